$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column "PES" before column B, shifting the existing
# B1:M1 header cells (and their formats) one column to the right,
# without touching the sheet's <cols> column-width definitions.
for ($col = 13; $col -ge 2; $col--) {
    $srcCell = $ws.Cells.Item(1, $col)
    $dstCell = $ws.Cells.Item(1, $col + 1)
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
    $dstCell.Value2 = $srcCell.Value2
}

$ws.Cells.Item(1, 2).Value2 = "PES"

[void]$ws.Range("B2").Select()
